# chore: Fix typos
# Corrects three German typos on the "Zuordnungstabelle" sheet:
#   - "lichtschranken" -> "Lichtschranken" (capitalised noun)
#   - "Lift fährt Rauf"   -> "Lift fährt rauf"
#   - "Lift fährt Runter" -> "Lift fährt runter"
# Also updates the sheet's current selection/print setup to match the
# state the workbook was left in after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zuordnungstabelle")

# --- Fix the typos -------------------------------------------------------
$ws.Range("D6").Value  = "Tür blockiert, Auslösung durch Lichtschranken"
$ws.Range("D18").Value = "Lift fährt rauf"
$ws.Range("D19").Value = "Lift fährt runter"

# --- Restore the sheet's used-range selection -----------------------------
$ws.Activate() | Out-Null
$ws.Range("B2:I19").Select() | Out-Null

# --- Page setup tweak recorded alongside the fix --------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
